$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / tab to reflect the new "through" date
$ws.Name = "Through 2021-09-27"

# Update the header label for the current (in-progress) month column
$ws.Range("B1").Value2 = "September 2021 (through September 27)"

# Updated / incremented counts (existing cells)
$ws.Range("K2").Value2 = 7
$ws.Range("K5").Value2 = 11
$ws.Range("AU10").Value2 = 5
$ws.Range("B11").Value2 = 6
$ws.Range("AL16").Value2 = 3
$ws.Range("BD20").Value2 = 4
$ws.Range("B22").Value2 = 2
$ws.Range("T22").Value2 = 4
$ws.Range("K43").Value2 = 4
$ws.Range("B56").Value2 = 2

# Newly populated cells (previously empty)
$ws.Range("T26").Value2 = 1
$ws.Range("AL54").Value2 = 1
$ws.Range("AU64").Value2 = 1
$ws.Range("K99").Value2 = 1
